$d = $word.ActiveDocument
$rng = $d.Range(0, 0)
$toc = $d.TablesOfContents.Add($rng, $true, 1, 7, $true, $null, $true, $true, $null, $true, $true, $true)
$d.Save()
